$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, [string]$text)
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '=T("' + $escaped + '")'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

Set-TextCell ($ws.Cells.Item(2, "D")) '64.016.65'
Set-TextCell ($ws.Cells.Item(2, "E")) '  +1.45%  '
Set-TextCell ($ws.Cells.Item(3, "D")) '3.093.05'
Set-TextCell ($ws.Cells.Item(3, "E")) '  +0.58%  '
Set-TextCell ($ws.Cells.Item(4, "E")) '  -0.09%  '
Set-TextCell ($ws.Cells.Item(5, "D")) '541.54'
Set-TextCell ($ws.Cells.Item(5, "E")) '  -0.73%  '
Set-TextCell ($ws.Cells.Item(6, "D")) '136.99'
Set-TextCell ($ws.Cells.Item(6, "E")) '  +0.70%  '
Set-TextCell ($ws.Cells.Item(7, "E")) '  -0.04%  '
Set-TextCell ($ws.Cells.Item(8, "D")) '3.088.63'
Set-TextCell ($ws.Cells.Item(8, "E")) '  +0.67%  '
Set-TextCell ($ws.Cells.Item(9, "E")) '  +1.33%  '
Set-TextCell ($ws.Cells.Item(10, "E")) '  +0.18%  '
Set-TextCell ($ws.Cells.Item(11, "D")) '6.44'
Set-TextCell ($ws.Cells.Item(11, "E")) '  -1.69%  '
Set-TextCell ($ws.Cells.Item(12, "E")) '  -0.21%  '
Set-TextCell ($ws.Cells.Item(13, "D")) '0.0000228'
Set-TextCell ($ws.Cells.Item(13, "E")) '  +5.72%  '
Set-TextCell ($ws.Cells.Item(14, "D")) '34.72'
Set-TextCell ($ws.Cells.Item(14, "E")) '  -0.28%  '
Set-TextCell ($ws.Cells.Item(15, "D")) '3.594.69'
Set-TextCell ($ws.Cells.Item(15, "E")) '  +0.41%  '
Set-TextCell ($ws.Cells.Item(16, "D")) '64.059.64'
Set-TextCell ($ws.Cells.Item(16, "E")) '  +1.34%  '
Set-TextCell ($ws.Cells.Item(17, "E")) '  +1.15%  '
Set-TextCell ($ws.Cells.Item(18, "D")) '3.090.85'
Set-TextCell ($ws.Cells.Item(18, "E")) '  +0.05%  '
Set-TextCell ($ws.Cells.Item(19, "E")) '  +1.11%  '
Set-TextCell ($ws.Cells.Item(20, "D")) '483.95'
Set-TextCell ($ws.Cells.Item(20, "E")) '  +0.25%  '
Set-TextCell ($ws.Cells.Item(21, "D")) '13.40'
Set-TextCell ($ws.Cells.Item(21, "E")) '  +0.34%  '
Set-TextCell ($ws.Cells.Item(22, "D")) '0.702'
Set-TextCell ($ws.Cells.Item(22, "E")) '  +0.53%  '
Set-TextCell ($ws.Cells.Item(23, "D")) '7.11'
Set-TextCell ($ws.Cells.Item(23, "E")) '  -0.41%  '
Set-TextCell ($ws.Cells.Item(24, "D")) '79.71'
Set-TextCell ($ws.Cells.Item(24, "E")) '  +2.78%  '
Set-TextCell ($ws.Cells.Item(25, "D")) '12.22'
Set-TextCell ($ws.Cells.Item(25, "E")) '  +0.76%  '
Set-TextCell ($ws.Cells.Item(26, "E")) '  +0.09%  '
Set-TextCell ($ws.Cells.Item(27, "E")) '  +0.01%  '
Set-TextCell ($ws.Cells.Item(28, "D")) '8.10'
Set-TextCell ($ws.Cells.Item(28, "E")) '  -1.30%  '
Set-TextCell ($ws.Cells.Item(29, "D")) '1.00'
Set-TextCell ($ws.Cells.Item(29, "E")) '  -0.14%  '
Set-TextCell ($ws.Cells.Item(30, "D")) '26.38'
Set-TextCell ($ws.Cells.Item(30, "E")) '  +0.42%  '
Set-TextCell ($ws.Cells.Item(31, "E")) '  -1.63%  '
Set-TextCell ($ws.Cells.Item(32, "D")) '1.15'
Set-TextCell ($ws.Cells.Item(32, "E")) '  +1.41%  '
Set-TextCell ($ws.Cells.Item(33, "D")) '57.85'
Set-TextCell ($ws.Cells.Item(33, "E")) '  -5.42%  '
Set-TextCell ($ws.Cells.Item(34, "E")) '  -5.32%  '
Set-TextCell ($ws.Cells.Item(35, "D")) '504.28'
Set-TextCell ($ws.Cells.Item(35, "E")) '  -4.52%  '
Set-TextCell ($ws.Cells.Item(36, "D")) '5.36'
Set-TextCell ($ws.Cells.Item(36, "E")) '  +4.47%  '
Set-TextCell ($ws.Cells.Item(37, "D")) '6.00'
Set-TextCell ($ws.Cells.Item(37, "E")) '  +1.53%  '
Set-TextCell ($ws.Cells.Item(38, "D")) '3.253.62'
Set-TextCell ($ws.Cells.Item(38, "E")) '  +5.28%  '
Set-TextCell ($ws.Cells.Item(39, "D")) '0.0399'
Set-TextCell ($ws.Cells.Item(39, "E")) '  +0.16%  '
Set-TextCell ($ws.Cells.Item(40, "D")) '0.0796'
Set-TextCell ($ws.Cells.Item(40, "E")) '  +1.16%  '
Set-TextCell ($ws.Cells.Item(41, "E")) '  +1.28%  '
Set-TextCell ($ws.Cells.Item(42, "B")) 'Cosmos'
Set-TextCell ($ws.Cells.Item(42, "C")) 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell ($ws.Cells.Item(42, "D")) '8.12'
Set-TextCell ($ws.Cells.Item(42, "E")) '  +0.78%  '
Set-TextCell ($ws.Cells.Item(43, "B")) 'dogwifhat'
Set-TextCell ($ws.Cells.Item(43, "C")) 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextCell ($ws.Cells.Item(43, "D")) '2.67'
Set-TextCell ($ws.Cells.Item(43, "E")) '  +0.93%  '
Set-TextCell ($ws.Cells.Item(44, "D")) '0.255'
Set-TextCell ($ws.Cells.Item(44, "E")) '  +1.18%  '
Set-TextCell ($ws.Cells.Item(45, "E")) '  +0.08%  '
Set-TextCell ($ws.Cells.Item(46, "B")) 'Monero'
Set-TextCell ($ws.Cells.Item(46, "C")) 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell ($ws.Cells.Item(46, "D")) '122.85'
Set-TextCell ($ws.Cells.Item(46, "E")) '  +1.16%  '
Set-TextCell ($ws.Cells.Item(47, "B")) 'Fetch.AI'
Set-TextCell ($ws.Cells.Item(47, "C")) 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextCell ($ws.Cells.Item(47, "D")) '2.05'
Set-TextCell ($ws.Cells.Item(47, "E")) '  +0.36%  '
Set-TextCell ($ws.Cells.Item(48, "D")) '0.0₃0530'
Set-TextCell ($ws.Cells.Item(48, "E")) '  +5.66%  '
Set-TextCell ($ws.Cells.Item(49, "D")) '24.66'
Set-TextCell ($ws.Cells.Item(49, "E")) '  +2.01%  '
Set-TextCell ($ws.Cells.Item(50, "E")) '  +2.03%  '
Set-TextCell ($ws.Cells.Item(51, "D")) '2.42'
Set-TextCell ($ws.Cells.Item(51, "E")) '  +2.20%  '
